$wb = $excel.ActiveWorkbook

# Rows (in column F, "想去人数") whose values changed, and their new values.
# The same update applies to both the "展览" sheet and the "全部类型" sheet,
# which carry duplicate data.
$updates = @{
    2  = 1587
    3  = 8938
    7  = 339
    8  = 167
    10 = 67
    11 = 3801
    13 = 377
    14 = 99
    15 = 4166
    19 = 4
    21 = 237
    23 = 2604
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Cells.Item($row, 6).Value = $updates[$row]
    }
}
